$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$shp = $m.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
Write-Host "tr.Text:" $tr.Text
$run1 = $tr.Runs(1, 1)
Write-Host "run1.Text: [" $run1.Text "]"
$runAll = $tr.Runs()
Write-Host "runAll.Count:" $runAll.Count
